$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "gamma01" block (rows 21-25): add two more predation-rate columns (p08, p09 exe HOME)
# Header row (row 22): new values in J:N
$ws.Range("J22").Value = 2.5
$ws.Range("K22").Value = 2.6
$ws.Range("L22").Value = 4.4
$ws.Range("M22").Value = 4.5
$ws.Range("N22").Value = 4.6

# Data row (row 23): new values in J:N
$ws.Range("J23").Value = 0.8143
$ws.Range("K23").Value = 0.8305
$ws.Range("L23").Value = 0.756
$ws.Range("M23").Value = 0.7525
$ws.Range("N23").Value = 0.7507

# Existing cell correction on row 24 (G24 value was a placeholder "1", now an actual result)
$ws.Range("G24").Value = 0.8403

# Move the active selection to K13 (also drops the stale topLeftCell scroll anchor)
$ws.Range("K13").Select()
